# Apply odds/score updates to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("BC4").Value = 126

# Row 5
$ws.Range("G5").Value = 1.95
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.33
$ws.Range("K5").Value = 1.95
$ws.Range("AA5").Value = 19
$ws.Range("AG5").Value = 9
$ws.Range("AO5").Value = 11
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 251

# Row 6
$ws.Range("G6").Value = 1.91
$ws.Range("H6").Value = 3.9
$ws.Range("I6").Value = 3.5
$ws.Range("Z6").Value = 17
$ws.Range("AA6").Value = 13
$ws.Range("AB6").Value = 19
$ws.Range("AC6").Value = 19
$ws.Range("AQ6").Value = 29
$ws.Range("AW6").Value = 6
$ws.Range("BA6").Value = 67

# Row 8
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65

# Row 9
$ws.Range("G9").Value = 1.57
$ws.Range("H9").Value = 3.8
$ws.Range("AX9").Value = 34

# Row 12
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11

# Row 13
$ws.Range("G13").Value = 3.8
$ws.Range("I13").Value = 2
$ws.Range("N13").Value = 7.5
$ws.Range("U13").Value = 2.05
$ws.Range("V13").Value = 1.7
$ws.Range("W13").Value = 9
$ws.Range("X13").Value = 19
$ws.Range("AC13").Value = 7.5
$ws.Range("AH13").Value = 8.5
$ws.Range("AJ13").Value = 17
$ws.Range("AM13").Value = 501
$ws.Range("AR13").Value = 126
$ws.Range("AU13").Value = 9
$ws.Range("AW13").Value = 3.75
$ws.Range("AX13").Value = 11

# Row 14
$ws.Range("G14").Value = 1.75
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 2.38
$ws.Range("Q14").Value = 1.95
$ws.Range("R14").Value = 1.95
$ws.Range("AK14").Value = 41

# Row 20
$ws.Range("G20").Value = 27
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = 3.55
$ws.Range("L20").Value = 1.29
$ws.Range("O20").Value = 1.07
$ws.Range("P20").Value = 6.7
$ws.Range("Q20").Value = 1.25
$ws.Range("R20").Value = 3.55
$ws.Range("S20").Value = 1.15
$ws.Range("T20").Value = 4.75
$ws.Range("U20").Value = 2.25
$ws.Range("V20").Value = 1.57
$ws.Range("X20").Value = 600
$ws.Range("Y20").Value = 120
$ws.Range("AB20").Value = 350
$ws.Range("AC20").Value = 26
$ws.Range("AE20").Value = 45
$ws.Range("AF20").Value = 175
$ws.Range("AG20").Value = 12.5
$ws.Range("AH20").Value = 7.4
$ws.Range("AI20").Value = 13
$ws.Range("AJ20").Value = 6.5
$ws.Range("AK20").Value = 11.25
$ws.Range("AL20").Value = 37
$ws.Range("AO20").Value = 200
$ws.Range("AP20").Value = 100
$ws.Range("AT20").Value = 4.75
$ws.Range("AU20").Value = 11.5
$ws.Range("AV20").Value = 80
$ws.Range("AW20").Value = 3.4
$ws.Range("AY20").Value = 13.5
$ws.Range("AZ20").Value = 7
$ws.Range("BA20").Value = 23
$ws.Range("BB20").Value = 150
